$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

$ws.Range("D5").Value = "not hispanic/latino = 40 (48.8%)`nwhite = 35 (42.7%)`nhispanic = 4 (4.9%)`npacific islander = 1 (1.2%)`nasian = 1 (1.2%)`nnative american = 1 (1.2%)"

$ws.Range("E5").Value = "not hispanic/latino = 205 (61.2%)`nwhite = 112 (33.4%)`nhispanic = 11 (3.3%)`nasian = 3 (0.9%)`nblack = 2 (0.6%)`npacific islander = 1 (0.3%)`nnative american = 1 (0.3%)"

$ws.Range("F5").Value = "not hispanic/latino = 104 (58.4%)`nwhite = 58 (32.6%)`nhispanic = 11 (6.2%)`nasian = 2 (1.1%)`npacific islander = 2 (1.1%)`nnative american = 1 (0.6%)"
